$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C2").Value = 123
$ws.Range("C3").Value = 123

$ws.Range("C19").Select()
